# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the per-fund holding detail for 2022-Q1.
# 2) Insert a new summary row at the top of the "总计" sheet's data
#    (right under the header) for 2022-Q1, pushing 2021-Q4 / 2021-Q3 down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper style constants (match look of existing header / index cells:
# bold font, thin box border, centered horizontally, top vertically).
# ---------------------------------------------------------------------
$xlCenter = -4108
$xlTop = -4160

function Format-HeaderLike($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = $xlCenter
    $rng.VerticalAlignment = $xlTop
}

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q1", inserted before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
Format-HeaderLike $q1.Range("B1:H1")

# Index column (A2:A9 = 0..7)
$colA = $q1.Range("A2:A9")
Format-HeaderLike $colA

# Fund detail rows
$funds = @(
    @("001725", "汇添富中国高端制造股票", "18.26", "90.98", "3.98", "0.7267", 10),
    @("015115", "汇添富中国高端制造股票D", "18.26", "90.98", "3.98", "0.7267", 10),
    @("007639", "汇添富3年封闭运作竞争优势灵活配置混合", "13.07", "60.50", "3.71", "0.4849", 4),
    @("002746", "汇添富多策略定期开放灵活配置混合", "4.55", "64.75", "4.14", "0.1884", 4),
    @("501063", "汇添富悦享定期开放混合", "2.18", "60.21", "3.71", "0.0809", 6),
    @("002567", "大成国家安全主题灵活配置混合", "0.34", "52.90", "4.16", "0.0141", 6),
    @("011761", "平安鑫瑞混合型证券投资基金A", "1.09", "20.41", "0.57", "0.0062", 6),
    @("011762", "平安鑫瑞混合型证券投资基金C", "0.34", "20.41", "0.57", "0.0019", 6)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $r = 2 + $i
    $row = $funds[$i]

    $q1.Cells.Item($r, 1).Value = $i

    Set-TextValue $q1.Cells.Item($r, 2) $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    Set-TextValue $q1.Cells.Item($r, 4) $row[2]
    Set-TextValue $q1.Cells.Item($r, 5) $row[3]
    Set-TextValue $q1.Cells.Item($r, 6) $row[4]
    Set-TextValue $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Update "总计" sheet: insert new top data row for 2022-Q1
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
Format-HeaderLike $total.Range("A2")
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 2.23

# Renumber the index column for the rows that got pushed down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
